# GradeBook.xlsx edit: add a new "Homework 4" column (H) with per-student
# grades/warm-ups, and update the selected cell to reflect where the user
# ended up after entering the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header for column H ---------------------------------------------
$ws.Range("H1").Value = "Homework 4"

# --- New "Homework 4" grades (column H), rows 2-15 ------------------------
# Formulas mirror the "points earned / points possible" style already used
# in columns D-G of this sheet.
$ws.Range("H2").Formula  = "=19/20"
$ws.Range("H3").Formula  = "=0"
$ws.Range("H4").Formula  = "=18/20"
$ws.Range("H5").Formula  = "=0"
$ws.Range("H6").Formula  = "=0"
$ws.Range("H7").Formula  = "=16/20"
$ws.Range("H8").Formula  = "=19/20"
$ws.Range("H9").Formula  = "=19/20"
$ws.Range("H10").Formula = "=0"
$ws.Range("H11").Formula = "=16/20"
$ws.Range("H12").Formula = "=0"
$ws.Range("H13").Formula = "=16/20"
$ws.Range("H14").Formula = "=17/20"
$ws.Range("H15").Formula = "=5/20"

# --- Column width for the new column H ------------------------------------
# (target OOXML stored width ~10.89 chars; ColumnWidth quantizes to pixel
# widths under the hood, so 10.0 is the closest achievable setting)
$ws.Columns.Item(8).ColumnWidth = 10.0

# --- Leave the selection where the editor ended up after data entry -------
$ws.Range("H9").Select() | Out-Null
